# Commit message: "Removed mass column requirement"
#
# The WiscSIMS column dictionary had a row defining a "Mass" column
# (ColNames="Mass", DictionaryColNames="Mass", unit="AMU", Type="Numeric").
# That row is no longer required, so delete it entirely - Excel shifts
# every row below it up by one, and the now-unreferenced "Mass"/"AMU"
# shared strings drop out of the saved workbook automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole row first (mirrors clicking the row header before
# deleting), then remove it.
$ws.Range("A19:XFD19").Select()
$ws.Rows.Item(19).Delete()
